$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# The sheet's rows 2 and 3 (the "Repayment" and "Waive interest" transaction
# rows) need their A:J content swapped, row 3's K:L formatting needs to pick
# up what row 4's K:L formatting used to be, and row 4's old (now unused)
# K:L cells need to disappear entirely.
#
# Plan: stash the two rows we still need (old row 2, old row 4) in a scratch
# area, delete row 2 outright (this shifts row 3 -> row 2 and row 4 -> row 3,
# which conveniently gives row 2 the right A:J content *and* gives row 3 the
# right K:L formatting, both carried over verbatim, since this engine moves
# cells - not styles - wholesale on a row delete/shift), then write the
# stashed rows back into row 3 (A:J only, keeping the inherited K:L) and
# row 4.
#
# NOTE: deleting row 2 with a shift shifts *every* row below it up by one -
# including the scratch rows - so the stashed data originally pasted at
# rows 100/101 ends up at rows 99/100 once the delete happens.

# Stash old row 2 (A:I) - values, then formats - at scratch row 100.
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A100:I100").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A100:I100").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Stash old row 4 (A:I) - values, then formats - at scratch row 101.
$ws.Range("A4:I4").Copy() | Out-Null
$ws.Range("A101:I101").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("A4:I4").Copy() | Out-Null
$ws.Range("A101:I101").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Drop row 2 entirely; row 3 -> row 2, row 4 -> row 3, row 4 becomes empty.
# Everything below also shifts up by one row (scratch rows 100/101 -> 99/100).
$ws.Range("A2:L2").Delete(-4162) | Out-Null              # xlShiftUp

# Put the old row 2 content (A:I) into row 3, on top of the inherited K:L.
$ws.Range("A99:I99").Copy() | Out-Null
$ws.Range("A3:I3").PasteSpecial(-4163) | Out-Null         # xlPasteValues
$ws.Range("A99:I99").Copy() | Out-Null
$ws.Range("A3:I3").PasteSpecial(-4122) | Out-Null         # xlPasteFormats

# Restore the old row 4 content (A:I only - row 4 no longer carries K:L
# cells at all, see below) into row 4.
$ws.Range("A100:I100").Copy() | Out-Null
$ws.Range("A4:I4").PasteSpecial(-4163) | Out-Null         # xlPasteValues
$ws.Range("A100:I100").Copy() | Out-Null
$ws.Range("A4:I4").PasteSpecial(-4122) | Out-Null         # xlPasteFormats

# Row 4's K:L cells are dropped entirely (they used to hold the now-moved
# style that row 3 inherited above).
$ws.Range("K4:L4").Clear() | Out-Null

# Clean up the scratch area.
$ws.Range("A99:I100").Clear() | Out-Null

# Column J (Loan Balance) does not follow the A:I swap - both rows 2 and 3
# settle on a value of 0, simply trading number-format styles (row 2 takes
# the "#,##0" style row 3 used to have, row 3 takes row 2's old General
# style).
$ws.Range("J2").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("J3").Copy() | Out-Null
$ws.Range("J200").PasteSpecial(-4122) | Out-Null   # stash J3's current (General) format
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null      # J3 <- J2's old "#,##0" format
$ws.Range("J200").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null      # J2 <- J3's old General format
$ws.Range("J200").Clear() | Out-Null

# Selection moves to J3.
$ws.Range("J3").Select() | Out-Null
